# Auto-generated Word COM-interop script
# Applies text replacements to bring before.docx in line with the target diff:
#  - update the header date
#  - update every division-fact cell in the table (cell-by-cell, in document order)

$d = $word.ActiveDocument

$replacements = @(
    @("2025-07-14 Monday", "2025-07-15 Tuesday"),
    @("44÷7=6, 2", "37÷4=9, 1"),
    @("95÷8=11, 7", "61÷9=6, 7"),
    @("24÷9=2, 6", "34÷8=4, 2"),
    @("25÷2=12, 1", "54÷8=6, 6"),
    @("61÷6=10, 1", "85÷8=10, 5"),
    @("73÷5=14, 3", "58÷3=19, 1"),
    @("92÷7=13, 1", "95÷2=47, 1"),
    @("76÷4=19, 0", "13÷2=6, 1"),
    @("73÷6=12, 1", "22÷5=4, 2"),
    @("96÷9=10, 6", "92÷8=11, 4"),
    @("69÷9=7, 6", "90÷2=45, 0"),
    @("20÷3=6, 2", "69÷9=7, 6"),
    @("89÷4=22, 1", "20÷4=5, 0"),
    @("70÷7=10, 0", "68÷4=17, 0"),
    @("56÷5=11, 1", "82÷5=16, 2"),
    @("45÷5=9, 0", "69÷6=11, 3"),
    @("32÷7=4, 4", "61÷9=6, 7"),
    @("46÷2=23, 0", "24÷4=6, 0"),
    @("34÷6=5, 4", "62÷8=7, 6"),
    @("44÷4=11, 0", "16÷6=2, 4"),
    @("68÷8=8, 4", "74÷8=9, 2"),
    @("75÷2=37, 1", "85÷5=17, 0"),
    @("81÷3=27, 0", "91÷7=13, 0"),
    @("64÷7=9, 1", "46÷3=15, 1"),
    @("72÷3=24, 0", "60÷3=20, 0"),
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                            $true, 1, $false, $new, 2)
}

